$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (from the commit diff).
$updates = [ordered]@{
    'D2' = '69.723.41'
    'E2' = '  -0.53%  '
    'D3' = '3.820.93'
    'E3' = '  +2.09%  '
    'E4' = '  +0.10%  '
    'D5' = '612.18'
    'E5' = '  -1.95%  '
    'D6' = '175.61'
    'E6' = '  -2.65%  '
    'D7' = '3.817.10'
    'E7' = '  +2.04%  '
    'D9' = '0.525'
    'E9' = '  -2.06%  '
    'E10' = '  -1.15%  '
    'D11' = '6.46'
    'E11' = '  +2.58%  '
    'D12' = '0.478'
    'E12' = '  -1.83%  '
    'D13' = '39.95'
    'E13' = '  -2.50%  '
    'E14' = '  -2.60%  '
    'D15' = '4.462.94'
    'E15' = '  +2.20%  '
    'D16' = '3.826.38'
    'E16' = '  +2.18%  '
    'D17' = '69.822.45'
    'E17' = '  -0.44%  '
    'D18' = '7.46'
    'E18' = '  -1.71%  '
    'E19' = '  -3.38%  '
    'D20' = '16.59'
    'E20' = '  -1.42%  '
    'D21' = '505.12'
    'E21' = '  -0.26%  '
    'D22' = '9.47'
    'E22' = '  +0.32%  '
    'D23' = '0.734'
    'E23' = '  +1.13%  '
    'B24' = 'Fetch.AI'
    'C24' = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
    'D24' = '2.45'
    'E24' = '  -3.70%  '
    'B25' = 'Litecoin'
    'C25' = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
    'D25' = '85.67'
    'E25' = '  -1.23%  '
    'E26' = '  +4.61%  '
    'D27' = '12.60'
    'E27' = '  -4.18%  '
    'E28' = '  -8.70%  '
    'E29' = '  +0.05%  '
    'D30' = '2.53'
    'E30' = '  +1.36%  '
    'D31' = '2.98'
    'E31' = '  +1.34%  '
    'D32' = '7.98'
    'E32' = '  +0.77%  '
    'D33' = '31.70'
    'E33' = '  +1.44%  '
    'E34' = '  -2.00%  '
    'D35' = '1.00'
    'E35' = '  +0.14%  '
    'E36' = '  -2.44%  '
    'E37' = '  -2.24%  '
    'E38' = '  +3.60%  '
    'D39' = '487.00'
    'E39' = '  +14.03%  '
    'E40' = '  -0.53%  '
    'B41' = 'Stacks'
    'C41' = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
    'D41' = '2.05'
    'E41' = '  -2.93%  '
    'B42' = 'dogwifhat'
    'C42' = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
    'D42' = '2.99'
    'E42' = '  +4.89%  '
    'D43' = '49.68'
    'E43' = '  -1.56%  '
    'D44' = '43.63'
    'E44' = '  -3.45%  '
    'E45' = '  -2.64%  '
    'D46' = '2.920.40'
    'E46' = '  -3.08%  '
    'D47' = '0.0361'
    'E47' = '  -1.10%  '
    'D48' = '139.81'
    'E48' = '  +1.89%  '
    'E49' = '  +0.03%  '
    'D50' = '26.83'
    'E50' = '  -2.15%  '
    'E51' = '  -4.45%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Force text storage so numeric-looking strings (e.g. "612.18",
    # "1.00", "0.0361") keep their exact original formatting instead
    # of being re-interpreted as numbers by Excel on assignment.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = "Normal"
}
